$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The LMIC column (C) previously stored a handful of rows as shared-string
# text "YES"/"NO" instead of real booleans like the rest of the column.
# Convert those remaining text cells to proper TRUE/FALSE boolean values so
# the whole column is consistent (this also lets Excel drop the now-unused
# "YES"/"NO" shared strings on save).
$ws.Range("C2").Value = $true    # Algeria -> LMIC = TRUE
$ws.Range("C5").Value = $false   # Bahrain -> LMIC = FALSE
$ws.Range("C19").Value = $true   # Egypt -> LMIC = TRUE
$ws.Range("C26").Value = $true   # Ghana -> LMIC = TRUE
$ws.Range("C30").Value = $false  # Hong Kong -> LMIC = FALSE
$ws.Range("C36").Value = $true   # Jamaica -> LMIC = TRUE
$ws.Range("C58").Value = $true   # Paraguay -> LMIC = TRUE
$ws.Range("C87").Value = $false  # United Arab Emirates -> LMIC = FALSE
$ws.Range("C91").Value = $false  # Zimbabwe -> LMIC = FALSE (was YES)
